# Append the new "out_vars" row for 2020-07-17 (SSA raw/clean data for July 17th).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force A48 to be stored as text (not auto-parsed into a date serial number),
# then restore the default "Normal" style so the cell carries no explicit
# style index (matching the other date cells in column A), leaving only the
# shared-string text value behind.
$ws.Range("A48").NumberFormat = "@"
$ws.Range("A48").Value = "2020-07-17"
$ws.Range("A48").Style = "Normal"

$ws.Range("B48").Value = 331298
$ws.Range("C48").Value = 382003
$ws.Range("D48").Value = 85877
$ws.Range("E48").Value = 38310
$ws.Range("F48").Value = 28.81
